$d = $word.ActiveDocument

# Locate the anchor point right after "...esteemed company" and before the
# period that starts "...The innovative culture...".
$rng = $d.Content
$rng.Find.Execute("esteemed company", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Insert the two new text pieces back-to-back, plainly, so the sentence
# reads correctly first. (The following ". The innovative culture..." text
# stays untouched in its original run.)
$p1start = $rng.End
$rng.InsertAfter(" for ")
$p1end = $rng.End

$rng.Collapse(0)
$p2start = $rng.End
$rng.InsertAfter("the Business/Data Analyst " + [char]0x2013 + " Finance MI & Analytics position")
$p2end = $rng.End

# Now force each of the two new pieces into its own run (distinct from its
# neighbors) by toggling a character property on and back off. Apply from
# last to first so earlier character offsets stay valid.
$r2 = $d.Range($p2start, $p2end)
$r2.Bold = 1
$r2.Bold = 0

$r1 = $d.Range($p1start, $p1end)
$r1.Bold = 1
$r1.Bold = 0
